# Applies the "Add files via upload" ranking update:
#   - C2 points value updated (18.01 -> 21.1)
#   - 58 new ranking rows appended (rows 158-215), including one
#     brand-new name in column C ("או ליהי ") on row 163
#   - selection/cursor moved to the new first-empty row (A216)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 21.1

$ws.Cells.Item(158, 1).Value = 'עדן ורד מרי'
$ws.Cells.Item(158, 2).Value = 1
$ws.Cells.Item(159, 1).Value = 'ליהי בראל'
$ws.Cells.Item(159, 2).Value = 1
$ws.Cells.Item(160, 1).Value = 'ירון גלפנד'
$ws.Cells.Item(160, 2).Value = 1
$ws.Cells.Item(161, 1).Value = 'אורי שטרנברג'
$ws.Cells.Item(161, 2).Value = 1
$ws.Cells.Item(162, 1).Value = 'איתי הראל'
$ws.Cells.Item(162, 2).Value = 1
$ws.Cells.Item(163, 1).Value = 'ירון גלפנד'
$ws.Cells.Item(163, 2).Value = 6
$ws.Cells.Item(163, 3).Value = 'או ליהי '
$ws.Cells.Item(164, 1).Value = 'ירון גלפנד'
$ws.Cells.Item(164, 2).Value = 6
$ws.Cells.Item(165, 1).Value = 'דן פימה'
$ws.Cells.Item(165, 2).Value = 1
$ws.Cells.Item(166, 1).Value = 'אביב ואסקז'
$ws.Cells.Item(166, 2).Value = 1
$ws.Cells.Item(167, 1).Value = 'אן מרש'
$ws.Cells.Item(167, 2).Value = 1
$ws.Cells.Item(168, 1).Value = 'יהלי דוייב'
$ws.Cells.Item(168, 2).Value = 1
$ws.Cells.Item(169, 1).Value = 'דפנה ברגשטיין'
$ws.Cells.Item(169, 2).Value = 1
$ws.Cells.Item(170, 1).Value = 'יולי יערי תליו'
$ws.Cells.Item(170, 2).Value = 1
$ws.Cells.Item(171, 1).Value = 'איתי בסטקר'
$ws.Cells.Item(171, 2).Value = 1
$ws.Cells.Item(172, 1).Value = 'יובל סטרוזר'
$ws.Cells.Item(172, 2).Value = 1
$ws.Cells.Item(173, 1).Value = 'תאיו ורד'
$ws.Cells.Item(173, 2).Value = 1
$ws.Cells.Item(174, 1).Value = 'תומר ששון'
$ws.Cells.Item(174, 2).Value = 1
$ws.Cells.Item(175, 1).Value = 'דן פימה'
$ws.Cells.Item(175, 2).Value = 6
$ws.Cells.Item(176, 1).Value = 'יהלי דוייב'
$ws.Cells.Item(176, 2).Value = 6
$ws.Cells.Item(177, 1).Value = 'הגר אגמון'
$ws.Cells.Item(177, 2).Value = 1
$ws.Cells.Item(178, 1).Value = 'יהלי דוייב'
$ws.Cells.Item(178, 2).Value = 1
$ws.Cells.Item(179, 1).Value = 'אן מרש'
$ws.Cells.Item(179, 2).Value = 1
$ws.Cells.Item(180, 1).Value = 'אורי שטרנברג'
$ws.Cells.Item(180, 2).Value = 1
$ws.Cells.Item(181, 1).Value = 'איתי הראל'
$ws.Cells.Item(181, 2).Value = 1
$ws.Cells.Item(182, 1).Value = 'יהלי דוייב'
$ws.Cells.Item(182, 2).Value = 6
$ws.Cells.Item(183, 1).Value = 'אן מרש'
$ws.Cells.Item(183, 2).Value = 6
$ws.Cells.Item(184, 1).Value = 'רומי הרשקוביץ'
$ws.Cells.Item(184, 2).Value = 1
$ws.Cells.Item(185, 1).Value = 'אביב ואסקז'
$ws.Cells.Item(185, 2).Value = 1
$ws.Cells.Item(186, 1).Value = 'ליהי בראל'
$ws.Cells.Item(186, 2).Value = 1
$ws.Cells.Item(187, 1).Value = 'יולי יערי תליו'
$ws.Cells.Item(187, 2).Value = 1
$ws.Cells.Item(188, 1).Value = 'תאיו ורד'
$ws.Cells.Item(188, 2).Value = 1
$ws.Cells.Item(189, 1).Value = 'ירון גלפנד'
$ws.Cells.Item(189, 2).Value = 1
$ws.Cells.Item(190, 1).Value = 'תומר ששון'
$ws.Cells.Item(190, 2).Value = 1
$ws.Cells.Item(191, 1).Value = 'יער אלביר'
$ws.Cells.Item(191, 2).Value = 1
$ws.Cells.Item(192, 1).Value = 'יהלי גודר'
$ws.Cells.Item(192, 2).Value = 1
$ws.Cells.Item(193, 1).Value = 'ליאם דיין '
$ws.Cells.Item(193, 2).Value = 1
$ws.Cells.Item(194, 1).Value = 'רומי הרשקוביץ'
$ws.Cells.Item(194, 2).Value = 6
$ws.Cells.Item(195, 1).Value = 'ליהי בראל'
$ws.Cells.Item(195, 2).Value = 6
$ws.Cells.Item(196, 1).Value = 'רומי הרשקוביץ'
$ws.Cells.Item(196, 2).Value = 1
$ws.Cells.Item(197, 1).Value = 'יובל סטרוזר'
$ws.Cells.Item(197, 2).Value = 1
$ws.Cells.Item(198, 1).Value = 'תאיו ורד'
$ws.Cells.Item(198, 2).Value = 1
$ws.Cells.Item(199, 1).Value = 'ליהי בראל'
$ws.Cells.Item(199, 2).Value = 1
$ws.Cells.Item(200, 1).Value = 'יהלי דוייב'
$ws.Cells.Item(200, 2).Value = 1
$ws.Cells.Item(201, 1).Value = 'תאיו ורד'
$ws.Cells.Item(201, 2).Value = 6
$ws.Cells.Item(202, 1).Value = 'יובל סטרוזר'
$ws.Cells.Item(202, 2).Value = 6
$ws.Cells.Item(203, 1).Value = 'אביב ואסקז'
$ws.Cells.Item(203, 2).Value = 1
$ws.Cells.Item(204, 1).Value = 'דן פימה'
$ws.Cells.Item(204, 2).Value = 1
$ws.Cells.Item(205, 1).Value = 'גלי זליג'
$ws.Cells.Item(205, 2).Value = 1
$ws.Cells.Item(206, 1).Value = 'דפנה ברגשטיין'
$ws.Cells.Item(206, 2).Value = 1
$ws.Cells.Item(207, 1).Value = 'יולי יערי תליו'
$ws.Cells.Item(207, 2).Value = 1
$ws.Cells.Item(208, 1).Value = 'אורי שטרנברג'
$ws.Cells.Item(208, 2).Value = 1
$ws.Cells.Item(209, 1).Value = 'איתי הראל'
$ws.Cells.Item(209, 2).Value = 1
$ws.Cells.Item(210, 1).Value = 'ליאם דיין '
$ws.Cells.Item(210, 2).Value = 1
$ws.Cells.Item(211, 1).Value = 'מעיין סטרוזר'
$ws.Cells.Item(211, 2).Value = 1
$ws.Cells.Item(212, 1).Value = 'איתי בסטקר'
$ws.Cells.Item(212, 2).Value = 1
$ws.Cells.Item(213, 1).Value = 'תומר ששון'
$ws.Cells.Item(213, 2).Value = 1
$ws.Cells.Item(214, 1).Value = 'גלי זליג'
$ws.Cells.Item(214, 2).Value = 6
$ws.Cells.Item(215, 1).Value = 'איתי בסטקר'
$ws.Cells.Item(215, 2).Value = 6

[void]$ws.Range("A216").Select()

